# Regenerate the "K" (strikeouts) column (column G) values for
# megill_trevor.xlsx save_data sheet, replacing the previous
# "Strike#" derived values with the recalculated K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K")
$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 2
    8  = 1
    9  = 0
    10 = 0
    11 = 2
    12 = 4
    13 = 3
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 2
    20 = 0
    21 = 2
    22 = 0
    23 = 0
    24 = 0
    25 = 1
    26 = 0
    27 = 0
    28 = 1
    29 = 1
    30 = 3
    31 = 2
    32 = 0
    33 = 1
    34 = 1
    35 = 1
    36 = 3
    37 = 4
    38 = 2
    39 = 1
    40 = 2
    41 = 0
    42 = 0
    43 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
